# Error Calculations and Plots
# Updates the "missing data" worksheet: fills in several previously-missing
# numeric values, blanks out several previously-populated ones, removes two
# rows (RM 232 and SC 92), and shifts the remaining bottom rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "RM 232" row (originally row 26) and, after the shift,
#     the "SC 92" row (originally row 28, now row 27) -------------------
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Value corrections for rows 3-25 (unaffected by the row deletions) --
$ws.Cells.Item(3, 3).Value = 11.2
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(5, 3).Value = ""
$ws.Cells.Item(5, 5).Value = ""
$ws.Cells.Item(9, 4).Value = -14.5
$ws.Cells.Item(10, 4).Value = -14.7
$ws.Cells.Item(11, 4).Value = -15.5
$ws.Cells.Item(12, 4).Value = -14.1
$ws.Cells.Item(12, 5).Value = ""
$ws.Cells.Item(13, 5).Value = -5.3
$ws.Cells.Item(14, 5).Value = -5.4
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = -7.3
$ws.Cells.Item(18, 4).Value = ""
$ws.Cells.Item(20, 4).Value = ""
$ws.Cells.Item(20, 5).Value = ""
$ws.Cells.Item(21, 3).Value = 12.7
$ws.Cells.Item(23, 3).Value = ""
$ws.Cells.Item(23, 5).Value = ""
$ws.Cells.Item(25, 5).Value = ""

# --- Value corrections for the shifted rows 26-33 -----------------------
$ws.Cells.Item(27, 5).Value = -10
$ws.Cells.Item(28, 5).Value = -5.9
$ws.Cells.Item(31, 4).Value = -13.7
$ws.Cells.Item(32, 3).Value = 10.5
$ws.Cells.Item(32, 4).Value = -14.7
$ws.Cells.Item(32, 5).Value = -6.4
